$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the publish date on the pub (row 3 -> JAMA Network Open entry):
# Aug 12, 2023 (45150) -> Aug 9, 2023 (45147)
$ws.Range("F3").Value = 45147

# Tweak the (wrapped-text) row height on the CV/row-2 entry back down to
# Excel's real max row height of 409.5pt.
$ws.Rows.Item(2).RowHeight = 409.5

# Leave the selection where the edit finished, on F4.
$ws.Range("F4").Select()
